# "Modif programmation staff CO 5"
#
# 1) Refresh the cached text of the auto-updating "datetimeFigureOut"
#    date field (12/05/2023 -> 17/05/2023) everywhere it is cached:
#    the slide master and every slide layout's date placeholder.
# 2) Reposition / resize a handful of the little signage picture
#    icons on slide 1.
# 3) Relocate the "Image 1" picture: remove it from its old spot and
#    put an equivalent picture ("Image 11") in its new spot, right
#    after "Image 30" (same picture, same size, just moved).

$EMU_PER_POINT = 12700

function ToPt($emu) {
    return $emu / $EMU_PER_POINT
}

function Get-ShapeByName($shapes, $name) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -eq $name) { return $sh }
    }
    return $null
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Date field re-cache: 12/05/2023 -> 17/05/2023
# ---------------------------------------------------------------
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDateShape = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) { $isDateShape = $true }
        } catch {
        }
        if ($isDateShape -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -like "*12/05/2023*") {
                $tr.Text = "17/05/2023"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DatePlaceholders $layout.Shapes
}

# ---------------------------------------------------------------
# 2) Reposition / resize picture shapes on slide 1
# ---------------------------------------------------------------
$s1 = $p.Slides.Item(1)

$img5 = Get-ShapeByName $s1.Shapes "Image 5"
$img5.Left = ToPt(-97382)
$img5.Top = ToPt(5841810)

$img6 = Get-ShapeByName $s1.Shapes "Image 6"
$img6.Left = ToPt(19723)
$img6.Top = ToPt(3558394)

$img12 = Get-ShapeByName $s1.Shapes "Image 12"
$img12.Left = ToPt(6888498)
$img12.Top = ToPt(5848146)
$img12.Width = ToPt(1041534)
$img12.Height = ToPt(969704)

$img7 = Get-ShapeByName $s1.Shapes "Image 7"
$img7.Left = ToPt(11364776)
$img7.Top = ToPt(3889800)

# ---------------------------------------------------------------
# 3) Move "Image 1" to a new spot (after "Image 30"), renaming it
#    "Image 11" in the process.
# ---------------------------------------------------------------
$img1 = Get-ShapeByName $s1.Shapes "Image 1"
$img1.Name = "Image 11"
$img1.Left = ToPt(1714027)
$img1.Top = ToPt(5922537)
